$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark column D as Text so that numeric-looking price strings
# (e.g. "1.009", "314.80") are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel's input parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.401.72'
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").Value = '1.874.27'
$ws.Range("E3").Value = '  -1.98%  '
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  -2.81%  '
$ws.Range("D5").Value = '314.80'
$ws.Range("E5").Value = '  -2.00%  '
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  -2.39%  '
$ws.Range("D7").Value = '0.5100'
$ws.Range("E7").Value = '  -2.24%  '
$ws.Range("D8").Value = '0.3928'
$ws.Range("E8").Value = '  -1.00%  '
$ws.Range("D9").Value = '0.08394'
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").Value = '1.107'
$ws.Range("E10").Value = '  -2.81%  '
$ws.Range("D11").Value = '6.244'
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").Value = '1.876.19'
$ws.Range("E12").Value = '  -2.29%  '
$ws.Range("D13").Value = '20.45'
$ws.Range("D14").Value = '7.249'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").Value = '1.009'
$ws.Range("E15").Value = '  -2.82%  '
$ws.Range("D16").Value = '0.00001104'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '90.88'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '0.06704'
$ws.Range("E18").Value = '  -2.07%  '
$ws.Range("D19").Value = '17.66'
$ws.Range("E19").Value = '  -2.17%  '
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("D21").Value = '5.941'
$ws.Range("E21").Value = '  -2.88%  '
$ws.Range("D22").Value = '28.445.26'
$ws.Range("E22").Value = '  -1.38%  '
$ws.Range("D23").Value = '11.10'
$ws.Range("E23").Value = '  -1.71%  '
$ws.Range("D24").Value = '2.256'
$ws.Range("E24").Value = '  -1.24%  '
$ws.Range("D25").Value = '2.092.90'
$ws.Range("E25").Value = '  -1.91%  '
$ws.Range("D26").Value = '160.99'
$ws.Range("E26").Value = '  -1.51%  '
$ws.Range("D27").Value = '20.64'
$ws.Range("E27").Value = '  -2.21%  '
$ws.Range("D28").Value = '2.370'
$ws.Range("E28").Value = '  -3.52%  '
$ws.Range("D29").Value = '126.17'
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("D31").Value = '1.048'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = '5.771'
$ws.Range("E32").Value = '  -3.94%  '
$ws.Range("D33").Value = '3.596'
$ws.Range("E33").Value = '  -2.68%  '
$ws.Range("D34").Value = '0.02430'
$ws.Range("E34").Value = '  -1.80%  '
$ws.Range("D35").Value = '0.06483'
$ws.Range("E35").Value = '  -2.74%  '
$ws.Range("D36").Value = '0.2180'
$ws.Range("E36").Value = '  -2.10%  '
$ws.Range("D37").Value = '8.876'
$ws.Range("E37").Value = '  -6.48%  '
$ws.Range("D38").Value = '1.261'
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("E39").Value = '  -0.70%  '
$ws.Range("D40").Value = '5.062'
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("D41").Value = '0.6420'
$ws.Range("E41").Value = '  -2.61%  '
$ws.Range("D42").Value = '11.14'
$ws.Range("E42").Value = '  -0.86%  '
$ws.Range("E43").Value = '  -2.55%  '
$ws.Range("D44").Value = '0.6046'
$ws.Range("E44").Value = '  -2.18%  '
$ws.Range("D45").Value = '13.08'
$ws.Range("E45").Value = '  -1.56%  '
$ws.Range("D46").Value = '3.694'
$ws.Range("E46").Value = '  -1.99%  '
$ws.Range("D47").Value = '2.005'
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("D48").Value = '121.92'
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("D49").Value = '1.205'
$ws.Range("E49").Value = '  -3.32%  '
$ws.Range("D50").Value = '1.190'
$ws.Range("E50").Value = '  -8.86%  '
$ws.Range("D51").Value = '0.06815'
$ws.Range("E51").Value = '  -2.44%  '

# Restore the default cell style on column D so no residual number-format
# override remains on the cells (keeps formatting identical to source).
$ws.Range("D2:D51").Style = "Normal"

